$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test entries (2 folders, 4 different animals) appended below the
# existing data, mirroring the existing Animal ID / Start Time / Length (min)
# layout and the time-of-day number format already used in column B.
$rows = @(
    @(121,  0.46574074074074073, 7),
    @(163,  0.46521990740740743, 10),
    @(164,  0.68958333333333333, 11),
    @(164,  0.69097222222222221, 12),
    @(1071, 0.69515046296296301, 6),
    @(1071, 0.69269675925925922, 11)
)

$startRow = 8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]

    $timeCell = $ws.Cells.Item($r, 2)
    $timeCell.Value = $data[1]
    $timeCell.NumberFormat = "h:mm:ss"

    $ws.Cells.Item($r, 3).Value = $data[2]
}

# Move the active selection past the newly added rows, as in the source file.
$ws.Range("C14").Select()

# The source workbook now carries an explicit (portrait) page setup.
$ws.PageSetup.Orientation = 1

